$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name order (swap pairs in the shared-string / display order) ---
# Row 18/19: Portugal overtakes Austria
$ws.Range("A18").Value = "Portugal"
$ws.Range("A19").Value = "Austria"

# Row 27/28: India overtakes Dinamarca
$ws.Range("A27").Value = "India"
$ws.Range("A28").Value = "Dinamarca"

# Row 123/124: Republica de Yibuti overtakes Brunei
$ws.Range("A123").Value = "Republica de Yibuti"
$ws.Range("A124").Value = "Brunei"

# --- Update timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 14:22"

# --- Update numeric data ---
# Row 14: Suiza
$ws.Range("E14").Value = 13227
$ws.Range("G14").Value = 37
$ws.Range("H14").Value = 858

# Row 15: Paises Bajos
$ws.Range("B15").Value = 20549
$ws.Range("C15").Value = 969
$ws.Range("E15").Value = 18051
$ws.Range("G15").Value = 147
$ws.Range("H15").Value = 2248

# Row 18: now Portugal
$ws.Range("B18").Value = 13141
$ws.Range("C18").Value = 699
$ws.Range("D18").Value = 196
$ws.Range("E18").Value = 12565
$ws.Range("F18").Value = 245
$ws.Range("G18").Value = 35
$ws.Range("H18").Value = 380

# Row 19: now Austria
$ws.Range("B19").Value = 12824
$ws.Range("C19").Value = 185
$ws.Range("D19").Value = 4512
$ws.Range("E19").Value = 8039
$ws.Range("F19").Value = 267
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 273

# Row 23: Suecia
$ws.Range("B23").Value = 8419
$ws.Range("C23").Value = 726
$ws.Range("E23").Value = 7527
$ws.Range("F23").Value = 678
$ws.Range("G23").Value = 96
$ws.Range("H23").Value = 687

# Row 27: now India
$ws.Range("B27").Value = 5480
$ws.Range("C27").Value = 129
$ws.Range("D27").Value = 468
$ws.Range("E27").Value = 4848
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 164

# Row 28: now Dinamarca
$ws.Range("B28").Value = 5386
$ws.Range("C28").Value = 315
$ws.Range("D28").Value = 1621
$ws.Range("E28").Value = 3547
$ws.Range("F28").Value = 127
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = 218

# Row 32: Rumania
$ws.Range("E32").Value = 4023
$ws.Range("G32").Value = 13
$ws.Range("H32").Value = 210

# Row 59
$ws.Range("B59").Value = 1343
$ws.Range("C59").Value = 61
$ws.Range("D59").Value = 179
$ws.Range("E59").Value = 1145
$ws.Range("F59").Value = 36
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 19

# Row 73
$ws.Range("B73").Value = 794
$ws.Range("C73").Value = 30
$ws.Range("E73").Value = 684

# Row 123: now Republica de Yibuti
$ws.Range("C123").Value = 45
$ws.Range("D123").Value = 25
$ws.Range("E123").Value = 110
$ws.Range("F123").Value = 0
$ws.Range("H123").Value = 0

# Row 124: now Brunei
$ws.Range("B124").Value = 135
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 91
$ws.Range("E124").Value = 43
$ws.Range("F124").Value = 3
$ws.Range("H124").Value = 1

# Row 162
$ws.Range("B162").Value = 21
$ws.Range("C162").Value = 1
$ws.Range("D162").Value = 2
